# Apply market-data refresh values (Sheets/Asura_Profits.xlsx diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 416.21622
$ws.Range("I19").Value = 386
$ws.Range("J19").Value = 439.2381
$ws.Range("K19").Value = 386
$ws.Range("L19").Value = 439.2381
$ws.Range("M19").Value = -211
$ws.Range("N19").Value = -789.2381
$ws.Range("H112").Value = 2533.75
$ws.Range("J112").Value = 2799.4285
$ws.Range("L112").Value = 8398.2855
$ws.Range("N112").Value = -10614.2855
$ws.Range("H137").Value = 2879.3225
$ws.Range("I137").Value = 1869.8
$ws.Range("J137").Value = 3360.0476
$ws.Range("K137").Value = 5609.4
$ws.Range("L137").Value = 10080.1428
$ws.Range("M137").Value = -3059.4
$ws.Range("N137").Value = -15180.1428
$ws.Range("H138").Value = 2967.01
$ws.Range("I138").Value = 1933.8462
$ws.Range("J138").Value = 3334.9863
$ws.Range("K138").Value = 5801.5386
$ws.Range("L138").Value = 10004.9589
$ws.Range("M138").Value = -661.5385999999999
$ws.Range("N138").Value = -20284.9589

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7483.976
$ws.Range("I32").Value = 6500.853
$ws.Range("J32").Value = 11662.25
$ws.Range("K32").Value = 6500.853
$ws.Range("L32").Value = 11662.25
$ws.Range("M32").Value = -6213.853
$ws.Range("N32").Value = -12236.25
$ws.Range("H61").Value = 3518.6667
$ws.Range("I61").Value = 3422.4
$ws.Range("K61").Value = 3422.4
$ws.Range("M61").Value = -3210.4
$ws.Range("H63").Value = 3438.2144
$ws.Range("I63").Value = 3012.2727
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 3012.2727
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2326.2727
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 3438.2144
$ws.Range("I66").Value = 3012.2727
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 15061.3635
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -11629.3635
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 1364.3396
$ws.Range("I74").Value = 1238.1428
$ws.Range("J74").Value = 1846.1818
$ws.Range("K74").Value = 1238.1428
$ws.Range("L74").Value = 1846.1818
$ws.Range("M74").Value = -364.1428000000001
$ws.Range("N74").Value = -3594.1818
$ws.Range("H77").Value = 1364.3396
$ws.Range("I77").Value = 1238.1428
$ws.Range("J77").Value = 1846.1818
$ws.Range("K77").Value = 6190.714
$ws.Range("L77").Value = 9230.909
$ws.Range("M77").Value = -1822.714
$ws.Range("N77").Value = -17966.909
$ws.Range("H122").Value = 6627.8
$ws.Range("I122").Value = 7683.9443
$ws.Range("J122").Value = 3912
$ws.Range("K122").Value = 23051.8329
$ws.Range("L122").Value = 11736
$ws.Range("M122").Value = -20601.8329
$ws.Range("N122").Value = -16636
$ws.Range("H132").Value = 5635.237
$ws.Range("I132").Value = 6005.0356
$ws.Range("K132").Value = 18015.1068
$ws.Range("M132").Value = -15485.1068
$ws.Range("H136").Value = 3518.6667
$ws.Range("I136").Value = 3422.4
$ws.Range("K136").Value = 10267.2
$ws.Range("M136").Value = -7717.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 77891.38
$ws.Range("I94").Value = 925.4286
$ws.Range("J94").Value = 167685
$ws.Range("K94").Value = 925.4286
$ws.Range("L94").Value = 167685
$ws.Range("M94").Value = -474.4286
$ws.Range("N94").Value = -168587
$ws.Range("H134").Value = 2399.6
$ws.Range("I134").Value = 2027.381
$ws.Range("J134").Value = 3268.111
$ws.Range("K134").Value = 6082.143
$ws.Range("L134").Value = 9804.332999999999
$ws.Range("M134").Value = -3547.143
$ws.Range("N134").Value = -14874.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4669.231
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 12
$ws.Range("N4").Value = -10224
$ws.Range("H31").Value = 1966.262
$ws.Range("I31").Value = 2139.3
$ws.Range("J31").Value = 1808.9546
$ws.Range("K31").Value = 2139.3
$ws.Range("L31").Value = 1808.9546
$ws.Range("M31").Value = -1844.3
$ws.Range("N31").Value = -2398.9546
$ws.Range("H34").Value = 1966.262
$ws.Range("I34").Value = 2139.3
$ws.Range("J34").Value = 1808.9546
$ws.Range("K34").Value = 2139.3
$ws.Range("L34").Value = 1808.9546
$ws.Range("M34").Value = -1937.3
$ws.Range("N34").Value = -2212.9546
$ws.Range("H58").Value = 2180880
$ws.Range("J58").Value = 1850.1666
$ws.Range("L58").Value = 1850.1666
$ws.Range("N58").Value = -2256.1666
$ws.Range("H107").Value = 1374
$ws.Range("I107").Value = 1374
$ws.Range("K107").Value = 1374
$ws.Range("M107").Value = 546
$ws.Range("H132").Value = 1505167.1
$ws.Range("I132").Value = 2705385
$ws.Range("K132").Value = 8116155
$ws.Range("M132").Value = -8113625
$ws.Range("H134").Value = 1827.7368
$ws.Range("J134").Value = 3204.6667
$ws.Range("L134").Value = 9614.000100000001
$ws.Range("N134").Value = -14684.0001
$ws.Range("H136").Value = 2180880
$ws.Range("J136").Value = 1850.1666
$ws.Range("L136").Value = 5550.4998
$ws.Range("N136").Value = -10650.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1263.9247
$ws.Range("J68").Value = 1581.7046
$ws.Range("L68").Value = 4745.1138
$ws.Range("N68").Value = -6367.1138
$ws.Range("H71").Value = 1263.9247
$ws.Range("J71").Value = 1581.7046
$ws.Range("L71").Value = 14235.3414
$ws.Range("N71").Value = -22347.3414
$ws.Range("H107").Value = 1275.4133
$ws.Range("J107").Value = 1671.4615
$ws.Range("L107").Value = 5014.3845
$ws.Range("N107").Value = -8854.3845
$ws.Range("H134").Value = 3083.16
$ws.Range("J134").Value = 4544.4443
$ws.Range("L134").Value = 13633.3329
$ws.Range("N134").Value = -23773.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6572
$ws.Range("I5").Value = 6572
$ws.Range("K5").Value = 6572
$ws.Range("M5").Value = -6460
$ws.Range("H70").Value = 440761.1
$ws.Range("I70").Value = 561445
$ws.Range("J70").Value = 6299
$ws.Range("K70").Value = 561445
$ws.Range("L70").Value = 6299
$ws.Range("M70").Value = -561175
$ws.Range("N70").Value = -6839
$ws.Range("H73").Value = 440761.1
$ws.Range("I73").Value = 561445
$ws.Range("J73").Value = 6299
$ws.Range("K73").Value = 561445
$ws.Range("L73").Value = 6299
$ws.Range("M73").Value = -560509
$ws.Range("N73").Value = -8171
$ws.Range("H80").Value = 3015.682
$ws.Range("I80").Value = 2814.7222
$ws.Range("J80").Value = 3920
$ws.Range("K80").Value = 2814.7222
$ws.Range("L80").Value = 3920
$ws.Range("M80").Value = -1816.7222
$ws.Range("N80").Value = -5916
$ws.Range("H83").Value = 3015.682
$ws.Range("I83").Value = 2814.7222
$ws.Range("J83").Value = 3920
$ws.Range("K83").Value = 14073.611
$ws.Range("L83").Value = 19600
$ws.Range("M83").Value = -9081.611000000001
$ws.Range("N83").Value = -29584
$ws.Range("H102").Value = 3899.8667
$ws.Range("I102").Value = 3899.923
$ws.Range("J102").Value = 3899.5
$ws.Range("K102").Value = 3899.923
$ws.Range("L102").Value = 3899.5
$ws.Range("M102").Value = -2277.923
$ws.Range("N102").Value = -7143.5
$ws.Range("H132").Value = 5332.6665
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2422.125
$ws.Range("I68").Value = 1704.6364
$ws.Range("K68").Value = 1704.6364
$ws.Range("M68").Value = -955.6364000000001
$ws.Range("H71").Value = 2422.125
$ws.Range("I71").Value = 1704.6364
$ws.Range("K71").Value = 8523.182000000001
$ws.Range("M71").Value = -4779.182000000001
$ws.Range("H132").Value = 5276.25
$ws.Range("I132").Value = 5382.423
$ws.Range("K132").Value = 16147.269
$ws.Range("M132").Value = -13617.269
$ws.Range("H136").Value = 3956
$ws.Range("J136").Value = 4485.7144
$ws.Range("L136").Value = 13457.1432
$ws.Range("N136").Value = -18557.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4029
$ws.Range("I132").Value = 4347.1665
$ws.Range("K132").Value = 13041.4995
$ws.Range("M132").Value = -10511.4995
$ws.Range("H136").Value = 3345.8076
$ws.Range("I136").Value = 3691.6155
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 11074.8465
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -8524.8465
$ws.Range("N136").Value = -14100
